# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values for the first data row
# (row 2) on both the zh-cn and de-de sheets with the freshly generated
# report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 12:54:33"
$wsZhCn.Range("H2").Value = "2016-03-22 12:54:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 12:54:37"
$wsDeDe.Range("H2").Value = "2016-03-22 12:55:02"
